$d = $word.ActiveDocument

$replacements = @(
    @{old="32÷6=5, 2"; new="52÷3=17, 1"},
    @{old="99÷2=49, 1"; new="80÷2=40, 0"},
    @{old="88÷6=14, 4"; new="78÷7=11, 1"},
    @{old="84÷7=12, 0"; new="74÷2=37, 0"},
    @{old="35÷6=5, 5"; new="54÷6=9, 0"},
    @{old="69÷7=9, 6"; new="72÷9=8, 0"},
    @{old="77÷4=19, 1"; new="75÷4=18, 3"},
    @{old="65÷4=16, 1"; new="64÷8=8, 0"},
    @{old="33÷8=4, 1"; new="82÷6=13, 4"},
    @{old="91÷7=13, 0"; new="36÷4=9, 0"},
    @{old="32÷3=10, 2"; new="17÷5=3, 2"},
    @{old="78÷5=15, 3"; new="44÷7=6, 2"},
    @{old="40÷6=6, 4"; new="48÷9=5, 3"},
    @{old="51÷8=6, 3"; new="28÷7=4, 0"},
    @{old="27÷3=9, 0"; new="99÷7=14, 1"},
    @{old="71÷4=17, 3"; new="85÷4=21, 1"},
    @{old="44÷6=7, 2"; new="20÷5=4, 0"},
    @{old="24÷3=8, 0"; new="11÷7=1, 4"},
    @{old="75÷6=12, 3"; new="74÷7=10, 4"},
    @{old="29÷5=5, 4"; new="54÷8=6, 6"},
    @{old="89÷4=22, 1"; new="62÷9=6, 8"},
    @{old="45÷2=22, 1"; new="15÷6=2, 3"},
    @{old="25÷2=12, 1"; new="58÷8=7, 2"},
    @{old="60÷3=20, 0"; new="32÷9=3, 5"},
    @{old="12÷4=3, 0"; new="26÷3=8, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
